$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# 1. The stray "_GoBack" bookmark currently sits right after the second
#    "Downloadable for free" run (end of the Hansen econometrics textbook
#    paragraph). Remove it from there ...
# -----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -----------------------------------------------------------------------
# 2. Remove the whole "Minimal requirement: " paragraph (a short heading
#    paragraph that duplicated/was folded into the following bullet).
# -----------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "Minimal requirement: `r") {
        $para.Range.Delete()
        break
    }
}

# -----------------------------------------------------------------------
# 3. ... and re-create "_GoBack" as an empty (collapsed) bookmark right at
#    the start of the "Undergraduate level ..." bullet paragraph.
# -----------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "Undergraduate level*") {
        $startRange = $para.Range.Duplicate
        $startRange.Collapse(1)
        $d.Bookmarks.Add("_GoBack", $startRange)
        break
    }
}

# -----------------------------------------------------------------------
# 4. Fix the "Updated" -> "Updates" typo in the GitHub-repo paragraph,
#    keeping the neighbouring runs ("U", "will be made as", ...) intact.
# -----------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*pdated will be made as*") {

        # 4a. Plain text fix (this will coalesce the surrounding runs into
        #     fewer runs -- that is fixed back up in 4b/4c below).
        $find = $para.Range.Duplicate
        $find.Find.Execute("pdated ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "pdates ", 2) | Out-Null

        # 4b. Re-establish the original run boundaries by toggling Bold
        #     on/off over each desired run span -- a formatting-only edit
        #     splits runs without re-merging unrelated text elsewhere.
        $anchor = $para.Range.Duplicate
        $anchor.Find.Execute("pdate", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0) | Out-Null
        $pdateStart = $anchor.Start
        $pdateEnd = $anchor.End

        function Reseat($rangeStart, $rangeEnd) {
            $r = $d.Range($rangeStart, $rangeEnd)
            $r.Bold = 1
            $r.Bold = 0
        }

        $openParenStart = $pdateStart - 2
        $openParenEnd   = $pdateStart - 1
        $uEnd           = $pdateStart
        $sEnd           = $pdateEnd + 1
        $spaceEnd       = $pdateEnd + 2

        Reseat $openParenStart $openParenEnd   # "("
        Reseat $openParenEnd   $uEnd           # "U"
        Reseat $pdateStart     $pdateEnd       # "pdate"
        Reseat $pdateEnd       $sEnd           # "s"
        Reseat $sEnd           $spaceEnd       # " "

        # 4c. Also keep "will be made as" separate from the trailing
        #     " the course progresses)" run, as in the original document.
        $tail = $d.Range($spaceEnd, $para.Range.End)
        $tail.Find.Execute("will be made as", $true, $false, $false, $false, `
                            $false, $true, 1, $false, "", 0) | Out-Null
        $tail.Bold = 1
        $tail.Bold = 0

        break
    }
}
